$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Propagate existing styles into the newly-used columns (G/H for most
# rows, H/I for row 2) by copying a same-row cell that already carries the
# right style, so styles.xml gains no new entries. ---

# Row 2 (merged title row, style 6) -> extend into H2:I2
$ws.Range("D2").Copy() | Out-Null
$ws.Range("H2:I2").PasteSpecial(-4122) | Out-Null

# Row 3 (header row, style 1) -> extend into G3:I3
$ws.Range("D3").Copy() | Out-Null
$ws.Range("G3:I3").PasteSpecial(-4122) | Out-Null

# Data rows 4,7,10 (style 1 across D:I) -> extend into G:I
$ws.Range("D4").Copy() | Out-Null
$ws.Range("G4:I4").PasteSpecial(-4122) | Out-Null

$ws.Range("D7").Copy() | Out-Null
$ws.Range("G7:I7").PasteSpecial(-4122) | Out-Null

$ws.Range("D10").Copy() | Out-Null
$ws.Range("G10:I10").PasteSpecial(-4122) | Out-Null

# Data rows 5,8,11 (style 1 across D:I) -> extend into G:I
$ws.Range("D5").Copy() | Out-Null
$ws.Range("G5:I5").PasteSpecial(-4122) | Out-Null

$ws.Range("D8").Copy() | Out-Null
$ws.Range("G8:I8").PasteSpecial(-4122) | Out-Null

$ws.Range("D11").Copy() | Out-Null
$ws.Range("G11:I11").PasteSpecial(-4122) | Out-Null

# Data rows 6,9,12 (style 1 across D:I) -> extend into G:I
$ws.Range("D6").Copy() | Out-Null
$ws.Range("G6:I6").PasteSpecial(-4122) | Out-Null

$ws.Range("D9").Copy() | Out-Null
$ws.Range("G9:I9").PasteSpecial(-4122) | Out-Null

$ws.Range("D12").Copy() | Out-Null
$ws.Range("G12:I12").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Title / header text ---
$ws.Range("D2").Value = "Time taken (ms)"

$ws.Range("B3").Value = "File name"
$ws.Range("C3").Value = "# Processes"
$ws.Range("D3").Value = "Test 1"
$ws.Range("E3").Value = "Test 2"
$ws.Range("F3").Value = "Test 3"
$ws.Range("G3").Value = "Test 4"
$ws.Range("H3").Value = "Test 5"
$ws.Range("I3").Value = "Average"

# --- File 1 block (rows 4-6) ---
$ws.Range("B4").Value = "File 1"
$ws.Range("C4").Value = "1 Process"
$ws.Range("D4").Value = 0.88
$ws.Range("E4").Value = 0.935
$ws.Range("F4").Value = 0.299
$ws.Range("G4").Value = 0.916
$ws.Range("H4").Value = 0.882
$ws.Range("I4").Formula = "=AVERAGE(D4:H4)"

$ws.Range("B5").Value = "1,000 lines"
$ws.Range("C5").Value = "2 Processes"
$ws.Range("D5").Value = 0.811
$ws.Range("E5").Value = 0.874
$ws.Range("F5").Value = 0.932
$ws.Range("G5").Value = 0.972
$ws.Range("H5").Value = 0.876
$ws.Range("I5").Formula = "=AVERAGE(D5:H5)"

$ws.Range("C6").Value = "4 Processes"
$ws.Range("D6").Value = 1.059
$ws.Range("E6").Value = 1.129
$ws.Range("F6").Value = 1.132
$ws.Range("G6").Value = 1.176
$ws.Range("H6").Value = 0.977
$ws.Range("I6").Formula = "=AVERAGE(D6:H6)"

# --- File 2 block (rows 7-9) ---
$ws.Range("B7").Value = "File 2"
$ws.Range("C7").Value = "1 Process"
$ws.Range("D7").Value = 3.397
$ws.Range("E7").Value = 3.394
$ws.Range("F7").Value = 3.44
$ws.Range("G7").Value = 3.346
$ws.Range("H7").Value = 3.357
$ws.Range("I7").Formula = "=AVERAGE(D7:H7)"

$ws.Range("B8").Value = "10,000 lines"
$ws.Range("C8").Value = "2 Processes"
$ws.Range("D8").Value = 2.175
$ws.Range("E8").Value = 2.215
$ws.Range("F8").Value = 2.061
$ws.Range("G8").Value = 2.137
$ws.Range("H8").Value = 2.22
$ws.Range("I8").Formula = "=AVERAGE(D8:H8)"

$ws.Range("C9").Value = "4 Processes"
$ws.Range("D9").Value = 1.775
$ws.Range("E9").Value = 1.63
$ws.Range("F9").Value = 1.715
$ws.Range("G9").Value = 1.678
$ws.Range("H9").Value = 1.625
$ws.Range("I9").Formula = "=AVERAGE(D9:H9)"

# --- File 3 block (rows 10-12) ---
$ws.Range("B10").Value = "File 3"
$ws.Range("C10").Value = "1 Process"
$ws.Range("D10").Value = 4.94
$ws.Range("E10").Value = 4.716
$ws.Range("F10").Value = 4.92
$ws.Range("G10").Value = 4.893
$ws.Range("H10").Value = 4.845
$ws.Range("I10").Formula = "=AVERAGE(D10:H10)"

$ws.Range("B11").Value = "100,000 lines"
$ws.Range("C11").Value = "2 Processes"
$ws.Range("D11").Value = 2.623
$ws.Range("E11").Value = 2.61
$ws.Range("F11").Value = 2.598
$ws.Range("G11").Value = 2.567
$ws.Range("H11").Value = 2.51
$ws.Range("I11").Formula = "=AVERAGE(D11:H11)"

$ws.Range("C12").Value = "4 Processes"
$ws.Range("D12").Value = 1.502
$ws.Range("E12").Value = 1.513
$ws.Range("F12").Value = 1.515
$ws.Range("G12").Value = 1.448
$ws.Range("H12").Value = 1.402
$ws.Range("I12").Formula = "=AVERAGE(D12:H12)"

# --- Extend the merged title cell to cover the new columns ---
$ws.Range("D2:G2").UnMerge() | Out-Null
$ws.Range("D2:I2").Merge() | Out-Null

# --- Misc view state to mirror the edit ---
$ws.Range("G16").Select()
